$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: bump the weekly report Volume/Number and date range ---
# A8 rich text: "Volume 32   Number  1" -> "Volume 32   Number  2"
$ws.Range("A8").Characters(21, 1).Text = "2"

# C9 rich text: "Report Covering the Week  12/30/2024  Through  1/5/2025"
#            -> "Report Covering the Week  1/6/2025  Through  1/12/2025"
# Replace the rightmost date run first so the first date'''s character
# offsets are unaffected by the differing replacement length.
$ws.Range("C9").Characters(48, 8).Text = "1/12/2025"
$ws.Range("C9").Characters(27, 10).Text = "1/6/2025"

# --- Crime-statistics grid (rows 15-33): refreshed weekly figures ---

# Simple numeric value updates (style unchanged)
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 4
$ws.Range("J16").Value = 4
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = -63.636363636363
$ws.Range("N16").Value = -90.243902439024
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -55.555555555555
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 40
$ws.Range("I17").Value = 10
$ws.Range("J17").Value = 11
$ws.Range("K17").Value = -9.090909090909
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = -33.333333333333
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 14
$ws.Range("H18").Value = -6.666666666666
$ws.Range("I18").Value = 6
$ws.Range("J18").Value = 6
$ws.Range("L18").Value = -33.333333333333
$ws.Range("M18").Value = -40
$ws.Range("N18").Value = -88.888888888888
$ws.Range("C19").Value = 7
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = 11.111111111111
$ws.Range("I19").Value = 10
$ws.Range("J19").Value = 12
$ws.Range("K19").Value = -16.666666666666
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -9.090909090909
$ws.Range("N19").Value = -54.545454545454
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = -24.137931034482
$ws.Range("I20").Value = 7
$ws.Range("J20").Value = 14
$ws.Range("K20").Value = -50
$ws.Range("L20").Value = 75
$ws.Range("M20").Value = -61.111111111111
$ws.Range("N20").Value = -93.518518518518
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -25
$ws.Range("F21").Value = 97
$ws.Range("G21").Value = 97
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 38
$ws.Range("J21").Value = 48
$ws.Range("K21").Value = -20.833333333333
$ws.Range("L21").Value = 35.714285714285
$ws.Range("M21").Value = -32.142857142857
$ws.Range("N21").Value = -84.232365145228
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = -27.272727272727
$ws.Range("F24").Value = 54
$ws.Range("G24").Value = 67
$ws.Range("H24").Value = -19.402985074626
$ws.Range("I24").Value = 19
$ws.Range("J24").Value = 25
$ws.Range("K24").Value = -24
$ws.Range("L24").Value = -26.923076923076
$ws.Range("M24").Value = -34.482758620689
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 200
$ws.Range("F25").Value = 9
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = -35.714285714285
$ws.Range("I25").Value = 4
$ws.Range("J25").Value = 7
$ws.Range("K25").Value = -42.857142857142
$ws.Range("L25").Value = -55.555555555555
$ws.Range("C26").Value = 6
$ws.Range("E26").Value = -25
$ws.Range("F26").Value = 31
$ws.Range("G26").Value = 39
$ws.Range("H26").Value = -20.512820512820
$ws.Range("I26").Value = 16
$ws.Range("J26").Value = 16
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 77.777777777777
$ws.Range("M26").Value = -38.461538461538
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 8
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 166.666666666667
$ws.Range("I28").Value = 5

# Cells changing from text (N/A style) to numeric
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J15").Value = 1
$ws.Range("J15").NumberFormat = '#,##0'
$ws.Range("K15").Value = 0
$ws.Range("K15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M15").Value = 0
$ws.Range("M15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N15").Value = 0
$ws.Range("N15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = -100
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J27").Value = 1
$ws.Range("J27").NumberFormat = '#,##0'
$ws.Range("K27").Value = 0
$ws.Range("K27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("D28").Value = 2
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("E28").Value = 0
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J28").Value = 2
$ws.Range("J28").NumberFormat = '#,##0'
$ws.Range("K28").Value = 150
$ws.Range("K28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L28").Value = 400
$ws.Range("L28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("D33").Value = 1
$ws.Range("D33").NumberFormat = '#,##0'
$ws.Range("E33").Value = -100
$ws.Range("E33").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J33").Value = 1
$ws.Range("J33").NumberFormat = '#,##0'
$ws.Range("K33").Value = -100
$ws.Range("K33").NumberFormat = '#,##0.0;"-"#,##0.0'

# Cells changing from numeric to text (N/A style)
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("C15").NumberFormat = "General"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("C27").NumberFormat = "General"
